$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update the "Periodo Mora" / "Valor Mora" data (rows 16-27) ---
# The periods are reordered from descending (2303..2204) to ascending (2204..2303),
# carrying each period's "Valor Mora" along with it.
$periods = @("2204","2205","2206","2207","2208","2209","2210","2211","2212","2301","2302","2303")
$valores = @(40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,40000,37333)

for ($i = 0; $i -lt 12; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
